# Update the 25 "two-digit divided by one-digit" problems in the single
# table on the page. The table has 20 rows x 5 columns, but only every
# 4th row (1, 5, 9, 13, 17) actually holds a problem; the rows in
# between are blank spacer rows.
#
# Several of the original problems are textually identical (e.g. "31÷5="
# and "19÷4=" each occur twice) but must become different new values
# depending on which cell they are in, so a document-wide Find/Replace
# can't be used safely. Instead we address each cell positionally via
# Table.Cell(row, column) and overwrite its Range.Text directly, which
# replaces only the digits/operator run text in that specific cell while
# preserving the existing run/paragraph formatting (font, size, etc.).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text  = "49÷5="
$t.Cell(1,2).Range.Text  = "87÷8="
$t.Cell(1,3).Range.Text  = "61÷5="
$t.Cell(1,4).Range.Text  = "37÷2="
$t.Cell(1,5).Range.Text  = "43÷4="

$t.Cell(5,1).Range.Text  = "19÷4="
$t.Cell(5,2).Range.Text  = "14÷3="
$t.Cell(5,3).Range.Text  = "12÷6="
$t.Cell(5,4).Range.Text  = "83÷9="
$t.Cell(5,5).Range.Text  = "16÷2="

$t.Cell(9,1).Range.Text  = "54÷7="
$t.Cell(9,2).Range.Text  = "89÷8="
$t.Cell(9,3).Range.Text  = "81÷4="
$t.Cell(9,4).Range.Text  = "19÷4="
$t.Cell(9,5).Range.Text  = "83÷6="

$t.Cell(13,1).Range.Text = "62÷4="
$t.Cell(13,2).Range.Text = "47÷8="
$t.Cell(13,3).Range.Text = "21÷8="
$t.Cell(13,4).Range.Text = "38÷9="
$t.Cell(13,5).Range.Text = "45÷9="

$t.Cell(17,1).Range.Text = "63÷4="
$t.Cell(17,2).Range.Text = "73÷9="
$t.Cell(17,3).Range.Text = "26÷9="
$t.Cell(17,4).Range.Text = "57÷6="
$t.Cell(17,5).Range.Text = "28÷4="
